$wb = $excel.ActiveWorkbook

# --- "Trends Status" sheet ---
$ws1 = $wb.Worksheets.Item("Trends Status")

$ws1.Range("C3").Value = 5
$ws1.Range("E3").Value = 18.5

$ws1.Range("C4").Value = 20
$ws1.Range("E4").Value = 74.09999999999999

$ws1.Range("E5").Value = 7.4

$ws1.Range("C7").Value = 25

# --- "Species qualification" sheet ---
$ws4 = $wb.Worksheets.Item("Species qualification")

$ws4.Range("C4").Value = 27
